$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row (columns A:T only, to keep the used-range/dimension
# tight as A1:T42) right before current row 7. This shifts the existing
# rows 7-41 down to 8-42.
$ws.Range("A7:T7").Insert()

# The newly inserted row 7 becomes a copy of row 6's data (the entry that
# used to be there before this week's new entry was recorded).
$ws.Range("A6:T6").Copy($ws.Range("A7:T7"))

# Row 6 now holds this week's new record: update the date and volume.
$ws.Range("D6").Value = 44670
$ws.Range("M6").Value = 200
